$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "finished date" column (C) for rows 6-8 with the same date
# already recorded in column B for each of those rows, re-using the
# existing date number format (style) from column B rather than minting a
# new one.
$rows = 6, 7, 8
foreach ($r in $rows) {
    $srcCell = $ws.Range("B$r")
    $dstCell = $ws.Range("C$r")

    # Copy the value first.
    $dstCell.Value2 = $srcCell.Value2

    # Then copy just the formatting (number format/style) from the source
    # cell so C gets the same date style already used in column B, instead
    # of Excel fabricating a brand-new custom number format.
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# The author's last recorded selection moved from C6 to C8.
$ws.Range("C8").Select() | Out-Null
